# Fill in the previously-empty log entry on row 27 of the Working Time Log
# table: a new Date / Start Time / End Time entry, with the "Work Time"
# (shared formula in column D) and the totals row (D36) recalculating
# automatically as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date column (A27) - 2017-11-20
$ws.Range("A27").Value = 43059

# Start Time (B27) - 16:05
$ws.Range("B27").Value = 0.67013888888888884

# End Time (C27) - 17:00
$ws.Range("C27").Value = 0.70833333333333337

# Recalculate the workbook so the shared "Work Time" formula in D27 and the
# Table1 totals-row formula in D36 pick up the new values.
$excel.Calculate()
